$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Sprint Completion Date: 07/18/2016 -> 07/17/2016
# ---------------------------------------------------------------------------
$pCompletion = $d.Paragraphs.Item(7).Range
$pCompletion.Find.Execute("07/18", $true, $false, $false, $false, $false, $true, 1, $false, "07/17", 2)

# ---------------------------------------------------------------------------
# 2) Revision Number: 1 -> 3
# ---------------------------------------------------------------------------
$pRevNum = $d.Paragraphs.Item(8).Range
$pRevNum.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "3", 2)

# ---------------------------------------------------------------------------
# 3) Revision Date: 07/17/2016 -> 07/18/2016
# ---------------------------------------------------------------------------
$pRevDate = $d.Paragraphs.Item(9).Range
$pRevDate.Find.Execute("07/17", $true, $false, $false, $false, $false, $true, 1, $false, "07/18", 2)

# ---------------------------------------------------------------------------
# 4) Insert three new Heading 1 sections ("actions to stop doing",
#    "Actions to start doing", "Actions to Keep doing") right before the
#    existing "Tasks Completed This Sprint" heading, each followed by a
#    blank (Normal) paragraph.
# ---------------------------------------------------------------------------
$tasksCompletedPara = $d.Paragraphs.Item(10)
$insertPoint = $d.Range($tasksCompletedPara.Range.Start, $tasksCompletedPara.Range.Start)
$newHeadingsBlock = "actions to stop doing`r`rActions to start doing`r`rActions to Keep doing`r`r"
$insertPoint.InsertBefore($newHeadingsBlock)

# The three blank spacer paragraphs inserted above inherited the Heading 1
# style from the paragraph mark they were inserted in front of; reset them
# back to Normal.
$d.Paragraphs.Item(11).Range.Style = "Normal"
$d.Paragraphs.Item(13).Range.Style = "Normal"
$d.Paragraphs.Item(15).Range.Style = "Normal"

# ---------------------------------------------------------------------------
# 5) The block of 9 identical blank "spacing" paragraphs that followed
#    "Tasks Completed This Sprint" must become:
#       blank, "Tasks not completed this Sprint", blank,
#       "Work Completion Rate ", blank (the last blank, untouched)
#    Delete the first 8 of the 9 blanks, keeping the 9th (last) one intact,
#    then insert the new blank/heading/blank/heading block in front of it.
# ---------------------------------------------------------------------------
$firstBlank = $d.Paragraphs.Item(17)
$eighthBlank = $d.Paragraphs.Item(24)
$d.Range($firstBlank.Range.Start, $eighthBlank.Range.End).Delete()

# Paragraph 17 is now the sole surviving (9th, originally-untouched) blank
# paragraph, immediately followed by "Images of prototype created".
$lastBlankStart = $d.Paragraphs.Item(17).Range.Start
$newBlock = "`rTasks not completed this Sprint`r`rWork Completion Rate `r"
$d.Range($lastBlankStart, $lastBlankStart).InsertBefore($newBlock)

# Tidy up paragraph styles: two brand-new blank paragraphs reset to Normal
# (clears any spacing formatting inherited from the split point), and the
# two new section headings set to Heading 1.
$d.Paragraphs.Item(17).Range.Style = "Normal"
$d.Paragraphs.Item(18).Range.Style = "Heading 1"
$d.Paragraphs.Item(19).Range.Style = "Normal"
$d.Paragraphs.Item(20).Range.Style = "Heading 1"
